# Corrección a Diebold Mariano: actualiza estadístico DM y p-value
# recalculados para cada par de N_Calib, y el resultado de
# "Significativo" (Mejor) de la fila N_Calib_1=20, N_Calib_2=200
# pasa de "Sí" a "No".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5807659188019426
$ws.Range("D2").Value = 0.5673005781462976
$ws.Range("C3").Value = 0.8954337011701105
$ws.Range("D3").Value = 0.3802483884933678
$ws.Range("C4").Value = 0.7155932952938638
$ws.Range("D4").Value = 0.4817690651435216
$ws.Range("C5").Value = 1.670560058280191
$ws.Range("D5").Value = 0.1089758933002558
$ws.Range("G5").Value = "No"
$ws.Range("C6").Value = 0.299283970472031
$ws.Range("D6").Value = 0.7675324718540257
$ws.Range("C7").Value = 0.2388927133560056
$ws.Range("D7").Value = 0.8134001248776319
$ws.Range("C8").Value = 1.471052863232385
$ws.Range("D8").Value = 0.1554366365665372
$ws.Range("C9").Value = -0.1120529584528923
$ws.Range("D9").Value = 0.9117975235621527
$ws.Range("C10").Value = 0.7614368176192444
$ws.Range("D10").Value = 0.4544816307227166
$ws.Range("C11").Value = 0.9514691085088003
$ws.Range("D11").Value = 0.3517033264307048
